# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status (col B) rows 2 & 3: "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Latest Target File (col E) rows 2 & 3: filled in with the same file the
#     row's "Source File Name" (col A) points to, as a hyperlink
#   - Latest Handback File (col F) rows 2 & 3: filled in with the same file the
#     row's "Latest Handoff File" (col C) points to, as a hyperlink
#   - Latest Handback DateTime (col G) rows 2 & 3: stamped with the handback time

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

function Get-HyperlinkForCell($ws, $cell) {
    $addr = $cell.Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl
        }
    }
    return $null
}

function Set-HyperlinkStyle($cell) {
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

$locales = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-02-17 10:17:31" },
    @{ Sheet = "de-de"; HandbackTime = "2016-02-17 10:17:52" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($row in 2, 3) {

        # Status -> Handed back
        $ws.Cells.Item($row, 2).Value = $newStatus

        $srcCell = $ws.Cells.Item($row, 1)       # A: Source File Name
        $handoffCell = $ws.Cells.Item($row, 3)   # C: Latest Handoff File
        $targetCell = $ws.Cells.Item($row, 5)    # E: Latest Target File
        $handbackCell = $ws.Cells.Item($row, 6)  # F: Latest Handback File

        $srcLink = Get-HyperlinkForCell $ws $srcCell
        $handoffLink = Get-HyperlinkForCell $ws $handoffCell

        # Latest Target File = same file as the source file, now confirmed in sync
        $targetCell.Value = $srcCell.Value2
        $ws.Hyperlinks.Add($targetCell, $srcLink.Address, [Type]::Missing, [Type]::Missing, $srcLink.TextToDisplay) | Out-Null
        Set-HyperlinkStyle $targetCell

        # Latest Handback File = same file as the latest handoff file
        $handbackCell.Value = $handoffCell.Value2
        $ws.Hyperlinks.Add($handbackCell, $handoffLink.Address, [Type]::Missing, [Type]::Missing, $handoffLink.TextToDisplay) | Out-Null
        Set-HyperlinkStyle $handbackCell

        # Latest Handback DateTime
        $ws.Cells.Item($row, 7).Value = $locale.HandbackTime
    }
}
